# Update the training data table and the active selection on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data corrections in row 2 (x_corrSteps, y_nrSteps, alienID)
$ws.Range("D2").Value = 7
$ws.Range("F2").Value = 3
$ws.Range("H2").Value = 46

# Move the active selection from D5 to C2
$ws.Range("C2").Select()
